$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.245.25'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.34%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.866.24'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.81%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.95'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4391'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3724'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07547'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9411'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.35%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.846.88'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.91%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.743'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.469'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06873'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '82.31'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009112'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.003'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.230.55'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.169'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.77'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.096.00'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.029'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.75'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.44'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.367'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.26'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.736'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09049'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8034'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -6.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.861'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.171'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.950'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.121'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05470'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01952'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.982'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +8.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.145'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5257'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1678'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.771'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.063'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.06772'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'PEPE'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.000002558'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.21%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4882'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.02%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '107.74'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.682'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.69%  '
